$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'35.398.43"
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = "'1.906.04"
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = "'244.41"
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('D6').Value = "'0.660"
$ws.Range('E6').Value = '  +6.10%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').Value = "'41.63"
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('E9').Value = '  +6.96%  '
$ws.Range('D10').Value = "'52.81"
$ws.Range('E10').Value = '  +12.82%  '
$ws.Range('D11').Value = "'0.0715"
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').Value = "'0.0998"
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('D13').Value = "'2.179.91"
$ws.Range('E13').Value = '  +2.59%  '
$ws.Range('D14').Value = "'12.05"
$ws.Range('E14').Value = '  +4.95%  '
$ws.Range('D15').Value = "'0.698"
$ws.Range('E15').Value = '  +3.05%  '
$ws.Range('D16').Value = "'1.909.52"
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').Value = "'4.86"
$ws.Range('E17').Value = '  +3.21%  '
$ws.Range('D18').Value = "'35.375.15"
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').Value = "'71.86"
$ws.Range('E19').Value = '  +2.90%  '
$ws.Range('D20').Value = "'0.0₃0822"
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('D21').Value = "'240.48"
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = "'12.53"
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('D23').Value = "'4.82"
$ws.Range('E23').Value = '  +2.16%  '
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('E25').Value = '  +25.64%  '
$ws.Range('D26').Value = "'2.30"
$ws.Range('E26').Value = '  +1.68%  '
$ws.Range('D27').Value = "'170.63"
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('D28').Value = "'8.44"
$ws.Range('E28').Value = '  +5.85%  '
$ws.Range('D29').Value = "'18.43"
$ws.Range('E29').Value = '  +4.49%  '
$ws.Range('E30').Value = '  +2.25%  '
$ws.Range('D31').Value = "'4.15"
$ws.Range('E31').Value = '  +3.94%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('B33').Value = 'BinanceUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D33').Value = "'1.01"
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'0.936"
$ws.Range('E34').Value = '  +11.99%  '
$ws.Range('E35').Value = '  +3.07%  '
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = "'2.03"
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'1.35"
$ws.Range('E38').Value = '  +3.84%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('D40').Value = "'0.0210"
$ws.Range('E40').Value = '  +4.54%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = "'0.0652"
$ws.Range('E41').Value = '  +17.22%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = "'16.36"
$ws.Range('E42').Value = '  +9.58%  '
$ws.Range('D43').Value = "'90.29"
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').Value = "'1.344.35"
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = "'2.40"
$ws.Range('E45').Value = '  +3.31%  '
$ws.Range('D46').Value = "'48.27"
$ws.Range('E46').Value = '  +38.93%  '
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('D48').Value = "'2.41"
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('D50').Value = "'2.090.23"
$ws.Range('E50').Value = '  +2.49%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.0697"
$ws.Range('E51').Value = '  +2.52%  '
